$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 106, shifting existing rows 106-119 down to 108-121
$ws.Rows("106:107").Insert()

# New row 106: "Especial" quality entry for the new price week
$ws.Range("A106").Value = 1
$ws.Range("B106").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C106").Value = "Arica y Parinacota"
$ws.Range("D106").Value = 44578
$ws.Range("E106").Value = 15
$ws.Range("F106").Value = "Fruta"
$ws.Range("G106").Value = 100108
$ws.Range("H106").Value = "Tropicales y subtropicales"
$ws.Range("I106").Value = 100108002
$ws.Range("J106").Value = "Mango"
$ws.Range("K106").Value = "Sin especificar"
$ws.Range("L106").Value = "Especial"
$ws.Range("M106").Value = 450
$ws.Range("N106").Value = 6500
$ws.Range("O106").Value = 7000
$ws.Range("P106").Value = 6750
$ws.Range("Q106").Value = "`$/bandeja 4 kilos"
$ws.Range("R106").Value = "Perú"
$ws.Range("S106").Value = 1688
$ws.Range("T106").Value = 4

# New row 107: "Primera" quality entry for the new price week
$ws.Range("A107").Value = 1
$ws.Range("B107").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C107").Value = "Arica y Parinacota"
$ws.Range("D107").Value = 44578
$ws.Range("E107").Value = 15
$ws.Range("F107").Value = "Fruta"
$ws.Range("G107").Value = 100108
$ws.Range("H107").Value = "Tropicales y subtropicales"
$ws.Range("I107").Value = 100108002
$ws.Range("J107").Value = "Mango"
$ws.Range("K107").Value = "Sin especificar"
$ws.Range("L107").Value = "Primera"
$ws.Range("M107").Value = 450
$ws.Range("N107").Value = 6500
$ws.Range("O107").Value = 7000
$ws.Range("P107").Value = 6750
$ws.Range("Q107").Value = "`$/bandeja 4 kilos"
$ws.Range("R107").Value = "Perú"
$ws.Range("S107").Value = 1688
$ws.Range("T107").Value = 4
